# Update weekly fruit/vegetable price records for rows 6, 7 and 8
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = 44516
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 900
$ws.Range("L6").Value = 1000
$ws.Range("M6").Value = 950
$ws.Range("P6").Value = 950

# Row 7
$ws.Range("D7").Value = 44505
$ws.Range("J7").Value = 440

# Row 8
$ws.Range("D8").Value = 44518
$ws.Range("J8").Value = 400
$ws.Range("K8").Value = 800
$ws.Range("L8").Value = 900
$ws.Range("M8").Value = 850
$ws.Range("P8").Value = 850
